$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.957.46"

$ws.Range("D3").Value = "2.467.64"
$ws.Range("E3").Value = "  -1.23%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'518.14"
$ws.Range("E5").Value = "  -3.58%  "

$ws.Range("D6").Value = "'130.81"
$ws.Range("E6").Value = "  -4.62%  "

$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("E8").Value = "  -2.46%  "

$ws.Range("D9").Value = "'0.0990"
$ws.Range("E9").Value = "  -2.18%  "

$ws.Range("E10").Value = "  -0.65%  "

$ws.Range("E11").Value = "  +0.07%  "

$ws.Range("D12").Value = "'0.342"
$ws.Range("E12").Value = "  -1.39%  "

$ws.Range("D13").Value = "2.905.49"
$ws.Range("E13").Value = "  -1.26%  "

$ws.Range("D14").Value = "57.876.21"
$ws.Range("E14").Value = "  -1.72%  "

$ws.Range("D15").Value = "'22.25"
$ws.Range("E15").Value = "  -3.36%  "

$ws.Range("E16").Value = "  -2.32%  "

$ws.Range("D17").Value = "2.468.17"
$ws.Range("E17").Value = "  -1.77%  "

$ws.Range("D18").Value = "'10.81"
$ws.Range("E18").Value = "  -2.87%  "

$ws.Range("E19").Value = "  -2.42%  "

$ws.Range("D20").Value = "'319.50"
$ws.Range("E20").Value = "  -1.32%  "

$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("D22").Value = "'5.73"
$ws.Range("E22").Value = "  -3.66%  "

$ws.Range("D23").Value = "'64.06"

$ws.Range("D24").Value = "'0.408"
$ws.Range("E24").Value = "  -3.13%  "

$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  +0.37%  "

$ws.Range("E26").Value = "  -3.47%  "

$ws.Range("D27").Value = "'7.32"
$ws.Range("E27").Value = "  -3.03%  "

$ws.Range("D28").Value = "0.0₃0751"
$ws.Range("E28").Value = "  -2.77%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'166.05"
$ws.Range("E29").Value = "  -1.02%  "

$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "'6.30"
$ws.Range("E30").Value = "  -5.85%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.69"
$ws.Range("E31").Value = "  -4.50%  "

$ws.Range("E32").Value = "  -1.93%  "

$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.24%  "

$ws.Range("D35").Value = "'18.04"
$ws.Range("E35").Value = "  -2.19%  "

$ws.Range("D36").Value = "'1.31"
$ws.Range("E36").Value = "  -9.87%  "

$ws.Range("E37").Value = "  -3.17%  "

$ws.Range("E38").Value = "  -4.57%  "

$ws.Range("E39").Value = "  -2.78%  "

$ws.Range("E40").Value = "  -4.37%  "

$ws.Range("D41").Value = "'272.23"
$ws.Range("E41").Value = "  -4.28%  "

$ws.Range("D42").Value = "'4.98"
$ws.Range("E42").Value = "  -3.22%  "

$ws.Range("D43").Value = "'0.591"
$ws.Range("E43").Value = "  -2.55%  "

$ws.Range("D44").Value = "'126.22"
$ws.Range("E44").Value = "  -4.95%  "

$ws.Range("D45").Value = "'0.0904"
$ws.Range("E45").Value = "  -2.41%  "

$ws.Range("E46").Value = "  -3.90%  "

$ws.Range("E47").Value = "  -3.28%  "

$ws.Range("D48").Value = "'17.07"
$ws.Range("E48").Value = "  -1.61%  "

$ws.Range("D49").Value = "1.732.33"
$ws.Range("E49").Value = "  -1.93%  "

$ws.Range("D50").Value = "'0.976"
$ws.Range("E50").Value = "  -1.00%  "

$ws.Range("D51").Value = "'4.69"
$ws.Range("E51").Value = "  -1.05%  "
